$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.Formula = '="58.643.25"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E2').Value = '  +1.53%  '
$r = $ws.Range('D3')
$r.Formula = '="3.158.71"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('E4').Value = '  -0.03%  '
$r = $ws.Range('D5')
$r.Formula = '="529.97"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E5').Value = '  -0.40%  '
$r = $ws.Range('D6')
$r.Formula = '="139.83"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E7').Value = '  -0.11%  '
$r = $ws.Range('D8')
$r.Formula = '="0.549"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E8').Value = '  +17.27%  '
$r = $ws.Range('D9')
$r.Formula = '="7.32"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E9').Value = '  +0.39%  '
$r = $ws.Range('D10')
$r.Formula = '="0.439"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E10').Value = '  +5.64%  '
$ws.Range('E11').Value = '  +4.34%  '
$ws.Range('E12').Value = '  +3.32%  '
$r = $ws.Range('D13')
$r.Formula = '="3.703.85"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E13').Value = '  +1.06%  '
$r = $ws.Range('D14')
$r.Formula = '="25.87"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E14').Value = '  +1.46%  '
$r = $ws.Range('D15')
$r.Formula = '="0.0000173"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E15').Value = '  +5.28%  '
$r = $ws.Range('D16')
$r.Formula = '="58.704.05"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E16').Value = '  +1.33%  '
$r = $ws.Range('D17')
$r.Formula = '="6.27"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E17').Value = '  +3.97%  '
$r = $ws.Range('D18')
$r.Formula = '="3.173.77"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E18').Value = '  +1.50%  '
$r = $ws.Range('D19')
$r.Formula = '="13.02"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E19').Value = '  +2.47%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$r = $ws.Range('D20')
$r.Formula = '="376.74"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E20').Value = '  +4.51%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$r = $ws.Range('D21')
$r.Formula = '="8.12"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E21').Value = '  +0.16%  '
$r = $ws.Range('D22')
$r.Formula = '="5.80"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E22').Value = '  +2.05%  '
$r = $ws.Range('D23')
$r.Formula = '="1.00"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  +5.16%  '
$r = $ws.Range('D25')
$r.Formula = '="69.82"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('E26').Value = '  +0.08%  '
$r = $ws.Range('D27')
$r.Formula = '="1.00"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E27').Value = '  +0.05%  '
$r = $ws.Range('D28')
$r.Formula = '="8.29"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E28').Value = '  +13.63%  '
$r = $ws.Range('D29')
$r.Formula = '="0.0₃0867"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E29').Value = '  -0.72%  '
$r = $ws.Range('D30')
$r.Formula = '="22.36"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E30').Value = '  +4.46%  '
$r = $ws.Range('D31')
$r.Formula = '="1.89"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E31').Value = '  +0.66%  '
$r = $ws.Range('D32')
$r.Formula = '="6.05"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E32').Value = '  -0.38%  '
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('E35').Value = '  +3.23%  '
$r = $ws.Range('D36')
$r.Formula = '="158.06"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('E37').Value = '  +5.15%  '
$r = $ws.Range('D38')
$r.Formula = '="24.97"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E38').Value = '  -3.44%  '
$ws.Range('E39').Value = '  +2.23%  '
$r = $ws.Range('D40')
$r.Formula = '="0.0693"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E40').Value = '  +2.99%  '
$r = $ws.Range('D41')
$r.Formula = '="2.649.78"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E41').Value = '  +6.15%  '
$ws.Range('E42').Value = '  +7.22%  '
$r = $ws.Range('D43')
$r.Formula = '="0.722"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E43').Value = '  +3.35%  '
$r = $ws.Range('D44')
$r.Formula = '="39.14"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E44').Value = '  +3.73%  '
$r = $ws.Range('D45')
$r.Formula = '="0.0289"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E45').Value = '  +7.69%  '
$r = $ws.Range('D46')
$r.Formula = '="1.00"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E46').Value = '  -0.03%  '
$r = $ws.Range('D47')
$r.Formula = '="3.199.73"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E47').Value = '  +0.92%  '
$r = $ws.Range('D48')
$r.Formula = '="0.105"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E48').Value = '  +15.03%  '
$ws.Range('E49').Value = '  +2.27%  '
$r = $ws.Range('D50')
$r.Formula = '="0.979"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E50').Value = '  -0.79%  '
$r = $ws.Range('D51')
$r.Formula = '="20.04"'
$r.Copy()
$r.PasteSpecial(-4163)
$ws.Range('E51').Value = '  +1.41%  '
